$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 6 (year 2025) metrics per commit "atualizei dados add e bibi"
$ws.Range("C6").Value = 405
$ws.Range("E6").Value = 98
$ws.Range("G6").Value = 24.19753086419753
$ws.Range("H6").Value = 75.80246913580247
